$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row with full words instead of abbreviations
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Move the active selection to A2, matching the saved view state
$ws.Range("A2").Select()
